$d = $word.ActiveDocument

# Replace in the main document body (addressee placeholder)
$d.Content.Find.Execute("QWREW", $true, $true, $false, $false, $false, $true, 1, $false, "QWR", 2)

# Replace in the header
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute("REW", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 2)
        $hdr.Range.Find.Execute("QWREW", $true, $true, $false, $false, $false, $true, 1, $false, "QWR", 2)
        $hdr.Range.Find.Execute("Rew", $true, $true, $false, $false, $false, $true, 1, $false, "Qwer", 2)
        $hdr.Range.Find.Execute("rew", $true, $true, $false, $false, $false, $true, 1, $false, "qwer", 2)
    }
}
